# Config.xlsx update: add System1 / SHA1Online assets, rename the
# logF_BusinessProcessName value from "Framework" to the new process
# name, and bump the MaxRetryNumber constant from 0 to 2.

$wb = $excel.ActiveWorkbook

# --- Settings sheet -------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# New asset rows describing System1 / SHA1Online endpoints and creds
$settings.Range("A6").Value = "System1_URL"
$settings.Range("A7").Value = "SHA1Online_URL"
$settings.Range("B6").Value = "https://acme-test.uipath.com/login"
$settings.Range("B7").Value = "https://crypt-online.ru/crypts/sha1/"
$settings.Range("A8").Value = "System1_Credential"
$settings.Range("B8").Value = "System1"

# Match the "General" number format Excel stamps onto freshly typed
# cells so B6/B7 pick up their own cellXfs entry (as in the authored diff).
$settings.Range("B6:B7").NumberFormat = "General"

# logF_BusinessProcessName value: Framework -> process name
$settings.Range("B5").Value = "Calculate Client Security Hash REFramework"

# --- Constants sheet --------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")
$constants.Range("B2").Value = 2
